# Update temperature problem boundary conditions on slides 1 and 2.
# "Wall temperature is 0°, except for a radiator at 100°"
#   -> "Wall temperature is 32°except for a radiator at 212°"

$p = $ppt.ActivePresentation

$degree = [char]0x00B0
$newText = "Wall temperature is 32" + $degree + "except for a radiator at 212" + $degree

for ($i = 1; $i -le 2; $i++) {
    $slide = $p.Slides.Item($i)
    $shape = $slide.Shapes.Item(2)
    $textRange = $shape.TextFrame.TextRange
    $paragraph = $textRange.Paragraphs(3)
    $run = $paragraph.Runs(1, 1)
    $run.Text = $newText
}
